$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "Docentes responsaveis" entry (Luiz Tadeu Fernandes Eleno) was added
# right below the existing one (Carlos Angelo Nunes, currently on row 14).
# Insert a fresh row at 14 which pushes the rest of the table
# (rows 14-24 -> 15-25) down by one, then populate it.
$ws.Rows.Item(14).Insert()

# Column B holds the plain value, column C mirrors it in red (the sheet's
# "changed data" column) - same pattern used by every other data row.
$ws.Cells.Item(14, 2).Value = "1176388 - Luiz Tadeu Fernandes Eleno"
$ws.Cells.Item(14, 3).Value = "1176388 - Luiz Tadeu Fernandes Eleno"

# Match the formatting used throughout the sheet for this kind of row
# (vertical-top, wrapped text; column C additionally in red) instead of
# whatever formatting Insert() may have copied down.
$ws.Cells.Item(14, 2).VerticalAlignment = -4160
$ws.Cells.Item(14, 2).WrapText = $true

$ws.Cells.Item(14, 3).VerticalAlignment = -4160
$ws.Cells.Item(14, 3).WrapText = $true
$ws.Cells.Item(14, 3).Font.Color = 255
